$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 = "Save", matching the style used by the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Data column H2:H14 -- the "Save" flag values
$saveValues = @(0, 1, 0, 1, 0, 0, 1, 0, 1, 0, 1, 0, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
